# LOB1235.xlsx update
# Inserts one new row after the "Objectives/Docentes responsaveis" block and
# re-populates the A/B/C "questionnaire" rows from row 13 downward with the
# corrected / extended content (new objectives text, new summary program,
# new full program, reshuffled method/criteria/recovery rows and new
# bibliography text), matching the upstream commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every row from 13 downward, down by one, inserting a fresh blank
# row 13. This naturally keeps all existing row heights / styles aligned
# with their new row numbers (verified to match the target row-height
# layout exactly), so nothing else needs to be re-stamped for height.
$ws.Rows.Item(13).Insert()

# ---- Row 10: Objetivos (PT) value, was showing the teacher's name by mistake ----
$ws.Range("B10").Value = "Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental."
$ws.Range("C10").Value = "Propiciar aos alunos conhecimento sobre os fundamentos, objetivos e métodos da Avaliação de Impacto Ambiental."

# ---- Row 13: Docentes responsaveis value (was implicitly missing a row) ----
$ws.Range("A13").Clear()

$ws.Range("B13").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160

$ws.Range("C13").Value = "5464150 - Mariana Consiglio Kasemodel"
$ws.Range("C13").Font.Color = 255
$ws.Range("C13").WrapText = $true
$ws.Range("C13").VerticalAlignment = -4160

# ---- Row 14: Programa resumido ----
$ws.Range("B14").Value = "Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais."
$ws.Range("C14").Value = "Os objetivos da Avaliação de Impacto Ambiental (AIA). O processo da AIA: triagem, definição do escopo, estudos de base, análise de impactos ambientais, mitigação, análise técnica e acompanhamento. Requisitos legais."

# ---- Row 15: Short syllabus (unchanged text, already correct after shift) ----

# ---- Row 16: Programa (full syllabus, PT) ----
$ws.Range("B16").Value = "Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso."
$ws.Range("C16").Value = "Conceitos básicos e definições. Origem e difusão da AIA. AIA e licenciamento: objetivos e fundamentos. Quadro legal e institucional brasileiro para a AIA. O processo de AIA e seus componentes. Etapas do planejamento e execução de um Estudo de Impacto Ambiental. Alternativas tecnológicas e de localização. Estudos de base e diagnóstico ambiental. Técnicas de identificação e previsão de impactos. Métodos e critérios de avaliação da importância dos impactos. Plano de gestão ambiental: medidas mitigadoras, compensatórias, de valorização e monitoramento. Tomada de decisão e acompanhamento. Estudos de caso."

# ---- Row 17: Syllabus (EN, unchanged text, already correct after shift) ----

# ---- Row 18: Avaliacao (label only, already correct after shift) ----

# ---- Row 19: Metodo ----
$ws.Range("B19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Range("C19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

# ---- Row 20: Criterio ----
$ws.Range("B20").Value = "Média ponderada de atividades e provas."
$ws.Range("C20").Value = "Média ponderada de atividades e provas."

# ---- Row 21: Norma de recuperacao ----
$ws.Range("B21").Value = "1 (uma) prova escrita"
$ws.Range("C21").Value = "1 (uma) prova escrita"

# ---- Row 22: Bibliografia ----
$ws.Range("B22").Value = "Bibliografia básicaSÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos. Oficina de textos: São Paulo, 2013. 583p.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão. Elsevier: Rio de Janeiro, 2019. 685p.Bibliografia complementar:COMPANHIA AMBIENTAL DO ESTADO DE SÃO PAULO - CETESB. MANUAL PARA ELABORAÇÃO DE ESTUDOS PARA O LICENCIAMENTO COM AVALIAÇÃO DE IMPACTO AMBIENTAL. Departamento de Desenvolvimento de Ações Estratégicas para o Licenciamento da Diretoria I- ID - CETESB. Anexo único, 2014. 250p."
$ws.Range("C22").Value = "Bibliografia básicaSÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos. Oficina de textos: São Paulo, 2013. 583p.CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão. Elsevier: Rio de Janeiro, 2019. 685p.Bibliografia complementar:COMPANHIA AMBIENTAL DO ESTADO DE SÃO PAULO - CETESB. MANUAL PARA ELABORAÇÃO DE ESTUDOS PARA O LICENCIAMENTO COM AVALIAÇÃO DE IMPACTO AMBIENTAL. Departamento de Desenvolvimento de Ações Estratégicas para o Licenciamento da Diretoria I- ID - CETESB. Anexo único, 2014. 250p."

# ---- Row 23: Requisitos (label only, already correct after shift) ----

# ---- Row 24: requirement detail row (already correct after shift) ----

Write-Host "LOB1235 sheet updated"
